$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")
    $val2 = $cellRow2.Value2
    $val3 = $cellRow3.Value2
    $cellRow2.Value = $val3
    $cellRow3.Value = $val2
}
